$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Zuletzt aktualisiert" SAVEDATE field result: 28.05.2014 12:52 -> 04.06.2014 10:27
# ------------------------------------------------------------------
$d.Content.Find.Execute("28.05.2014 12:52", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "04.06.2014 10:27", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Locate the two paragraphs at the end of the 04.06.2014 minutes section:
#      - "Teilnehmer: awe, bsh, cgu" (holds the _GoBack bookmark at its end)
#      - the empty placeholder paragraph right after the "Wir machen..." bullet
# ------------------------------------------------------------------
$teilnehmerPara = $null
$placeholderPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "Teilnehmer: awe, bsh, cgu*") {
        $teilnehmerPara = $i
    }
    if ($teilnehmerPara -ne $null -and $i -gt $teilnehmerPara -and $t -like "Wir machen z.Z.*") {
        $placeholderPara = $i + 1
        break
    }
}

# ------------------------------------------------------------------
# 3) Strip the _GoBack bookmark from the "Teilnehmer: awe, bsh, cgu" paragraph
#    (it moves further down into the new "Awe: ..." bullet).
# ------------------------------------------------------------------
$tp = $d.Paragraphs.Item($teilnehmerPara)
$tpXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:pStyle w:val="berschrift3Zwischentitel"/></w:pPr>
  <w:r><w:t>Teilnehmer: awe, bsh, cgu</w:t></w:r>
  <w:r><w:br/></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$tp.Range.InsertXML($tpXml) | Out-Null

# ------------------------------------------------------------------
# 4) Replace the trailing empty placeholder paragraph with the new minutes:
#      - Heading2 "31.07.2014"
#      - Zwischentitel "Ressourcenplanung"
#      - Aufzaehlung bullet "Awe: ..." (re-homes the _GoBack bookmark)
# ------------------------------------------------------------------
$pp = $d.Paragraphs.Item($placeholderPara)
$ppXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:pStyle w:val="Heading2"/></w:pPr>
  <w:r><w:t>31</w:t></w:r>
  <w:r><w:t>.0</w:t></w:r>
  <w:r><w:t>7</w:t></w:r>
  <w:r><w:t>.2014</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="berschrift3Zwischentitel"/></w:pPr>
  <w:r><w:t>Ressourcenplanung</w:t></w:r>
  <w:r><w:br/></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="Aufzhlung"/></w:pPr>
  <w:r><w:t xml:space="preserve">Awe: </w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t>Beat wird vom Html UI Projekt abgezogen und wird für zwei Monate bei Gothaer arbeiten. Samuel Moser wurde darüber informiert, dass sich darum die Lieferung vom Html UI um ca. 1 Monat verzögern wird. Eine Ersatzperson macht für den Zeitraum von zwei Monaten aufgrund der Einarbeitungszeit wenig Sinn.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$pp.Range.InsertXML($ppXml) | Out-Null

Write-Output "edit complete"
